$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# A handback just completed for the de-de file (and zh-cn is back in sync
# too), so the generated report is refreshed: the "Status" column now reads
# "Handed back: in sync with en-US" instead of "Ready for handoff", the
# "Latest Handback DateTime" timestamps move forward, and the stale
# "Error Detail" text (about the handback not being the latest version) is
# cleared out now that the handback is current.
# ---------------------------------------------------------------------------

$newStatus = "Handed back: in sync with en-US"

# Overview sheet: per-language status columns (zh-cn / de-de)
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# zh-cn detail sheet
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-08-30 06:52:28"
# Clear the "Error Detail" cell back to an empty string (keep the cell / shared
# string, don't just delete it - force text type via a leading quote then
# strip the resulting quote-prefix style so it serializes as plain s="0").
$wsZhCn.Range("P2").Value = "'"
$wsZhCn.Range("P2").Style = "Normal"

# de-de detail sheet
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-08-30 06:52:35"
$wsDeDe.Range("P2").Value = "'"
$wsDeDe.Range("P2").Style = "Normal"

# ---------------------------------------------------------------------------
# Column widths: the Status columns grew (longer status text) and the now-
# empty Error Detail columns shrank. Re-fit the columns to their new content
# (target widths measured from the regenerated report: ~29.98 chars for the
# Status columns, ~13.75 chars for the now-empty Error Detail columns).
# ---------------------------------------------------------------------------
$statusColWidth = 29.166666666666668
$errorDetailColWidth = 12.833333333333334

$wsOverview.Columns.Item(5).ColumnWidth = $statusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $statusColWidth
$wsZhCn.Columns.Item(16).ColumnWidth = $errorDetailColWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $statusColWidth
$wsDeDe.Columns.Item(16).ColumnWidth = $errorDetailColWidth
